# Update page-view counts on the "展览" (sheet1) and "全部类型" (sheet4) sheets.
# Note: row 9's value (F9) diverges by 1 between the two sheets in the source data,
# so each sheet is updated with its own explicit value.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 684
$wsExhibit.Range("F3").Value = 526
$wsExhibit.Range("F4").Value = 39
$wsExhibit.Range("F5").Value = 23
$wsExhibit.Range("F6").Value = 53
$wsExhibit.Range("F8").Value = 3329
$wsExhibit.Range("F9").Value = 4266
$wsExhibit.Range("F10").Value = 118

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 684
$wsAll.Range("F3").Value = 526
$wsAll.Range("F4").Value = 39
$wsAll.Range("F5").Value = 23
$wsAll.Range("F6").Value = 53
$wsAll.Range("F8").Value = 3329
$wsAll.Range("F9").Value = 4267
$wsAll.Range("F10").Value = 118
